$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Cronograma")

# Sprint 2 section moved from "En Progreso"/"Sin empezar" to "Finalizado",
# except the Sprint 3 header and "Función de busqueda" which moved to "En Progreso".
$ws.Range("F10").Value = "Finalizado"
$ws.Range("F11").Value = "Finalizado"
$ws.Range("F12").Value = "Finalizado"
$ws.Range("F13").Value = "Finalizado"
$ws.Range("F14").Value = "Finalizado"
$ws.Range("F15").Value = "Finalizado"
$ws.Range("F16").Value = "En Progreso"
$ws.Range("F17").Value = "Finalizado"
$ws.Range("F18").Value = "En Progreso"

# Update the current selection to match the saved cursor position.
$ws.Range("G17").Select()
